$d = $word.ActiveDocument
$BRK = [char]11

function Replace-InRange($range, $old, $new) {
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output ("WARNING: replace failed for: " + $old)
    }
    return $ok
}

function Replace-Doc($old, $new) {
    Replace-InRange $d.Content $old $new
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Doc "Unraveling the Enigma of Dark Matter" "Delving into the Realm of Mathematics: The Symphony of Numbers"

# ---------------------------------------------------------------------------
# Byline (author name). Target renders as "Albert D. Richards".
# ---------------------------------------------------------------------------
Replace-Doc "Alexej Steinn" "Albert D. Richards"

# ---------------------------------------------------------------------------
# Email address paragraph: "alexsteinn@constellation" + "." + "edu"
#                       -> "richardsalbert@eduworld" + "." + "in"
# Scope to paragraph 3 and do "edu"->"in" first is unnecessary since we scope
# narrowly, but we still replace the prefix before introducing new "edu" text.
# ---------------------------------------------------------------------------
$emailPara = $d.Paragraphs(3).Range
Replace-InRange $emailPara "edu" "in"
Replace-InRange $d.Paragraphs(3).Range "alexsteinn@constellation" "richardsalbert@eduworld"

# ---------------------------------------------------------------------------
# Body paragraph (paragraph 5) - sentence by sentence, mirroring the source
# edit. This paragraph mixes dark-matter content with new mathematics content.
# ---------------------------------------------------------------------------
$bodyPara = $d.Paragraphs(5).Range

Replace-InRange $d.Paragraphs(5).Range `
  "In the vast expanse of the cosmos, concealed within the depths of galaxies, lies a mysterious entity: dark matter" `
  "In the realm of academics, few disciplines capture the beauty of abstract thought and intellectual rigor like mathematics"

Replace-InRange $d.Paragraphs(5).Range `
  " This elusive substance, invisible to our eyes and undetected by traditional telescopes, exerts a gravitational influence that shapes the universe" `
  " It is a science that surrounds us, permeating every aspect of our world, from the delicate petals of flowers to the vast expanses of the cosmos"

# Remove " Its existence...reality" sentence and its trailing period entirely.
Replace-InRange $d.Paragraphs(5).Range `
  " Its existence, though intangible, is inferred through its gravitational effects on visible matter, providing a glimpse into the enigmas that permeate the fabric of reality." `
  ""

Replace-InRange $d.Paragraphs(5).Range `
  " Embarking on a voyage of discovery, scientists endeavor to unravel the enigma of dark matter, seeking to illuminate its properties, unravel its enigmatic nature, and decipher its intricate interplay with the universe" `
  " As high school students embark on their mathematical journey, they discover that mathematics is not merely a collection of formulas and equations, but rather an intricate symphony of numbers, patterns, and relationships"

Replace-InRange $d.Paragraphs(5).Range `
  "Traveling across galaxies, astronomers observe the peculiar motions of stars within clusters, revealing a disparity between the expected gravitational forces and the observed dynamics" `
  "Mathematics is the language of logic and reason, a tool that empowers us to understand the world around us"

Replace-InRange $d.Paragraphs(5).Range `
  " This discrepancy suggests the presence of unseen mass, exerting a gravitational pull on visible matter" `
  " The study of mathematical concepts offers a rigorous framework for problem-solving and a systematic approach to analyzing complex scenarios"

Replace-InRange $d.Paragraphs(5).Range `
  " This invisible entity, dubbed dark matter, constitutes approximately 85% of the matter in the universe, dwarfing the contribution of the familiar matter that surrounds us" `
  " Its beauty lies in its ability to transcend cultural and linguistic barriers, uniting humanity in a shared understanding of the universe's fundamental laws"

Replace-InRange $d.Paragraphs(5).Range `
  "Delving further into the mysteries of dark matter, physicists have proposed numerous theories attempting to explain its enigmatic nature" `
  "As students immerse themselves in the world of mathematics, they develop a profound sense of curiosity and a keen eye for patterns"

Replace-InRange $d.Paragraphs(5).Range `
  " Among them, the Weakly Interacting Massive Particle (WIMP) hypothesis has gained significant attention" `
  " Mathematical problems often require a creative approach, challenging conventional thinking and encouraging students to explore unconventional solutions"

Replace-InRange $d.Paragraphs(5).Range `
  " WIMPs are hypothetical particles predicted by certain extensions of the Standard Model of particle physics, possessing weak interactions with ordinary matter" `
  " It fosters an environment where intellectual exploration and perseverance are celebrated, helping students overcome challenges with resilience and determination"

# Final sentence of the old body becomes a plain line break, followed by a
# large new block of inserted content (the final "." run right after it is
# left completely untouched, so don't swallow it here).
$newTail = $BRK + "Body:" + $BRK + $BRK + "* Paragraph 1: The Aesthetic Charm of Mathematics:" + $BRK + $BRK + `
  "Mathematics has an inherent aesthetic appeal that captivates the human mind." + `
  " The harmony and balance found in mathematical patterns, equations, and geometric constructions have inspired artists, composers, and scientists throughout history." + `
  " This aesthetic dimension of mathematics not only enhances its beauty but also reinforces its universal nature." + `
  " It has been said that mathematics is the music of reason, with its own unique melodies, rhythms, and harmonies." + `
  $BRK + $BRK + "* Paragraph 2: Mathematics in Everyday Life:" + $BRK + $BRK + `
  "The practical applications of mathematics extend far beyond the classroom walls." + `
  " From the design of bridges and buildings to the intricacies of finance and economics, mathematics plays a vital role in shaping our modern world." + `
  " Students discover how mathematics helps predict weather patterns, model epidemics, and optimize transportation systems." + `
  " They learn the language of data, developing invaluable skills in analyzing and interpreting information, preparing them for careers in an increasingly data-driven society." + `
  $BRK + $BRK + "* Paragraph 3: Mathematical Thinking and Personal Growth:" + $BRK + $BRK + `
  "Mathematics is not solely about solving equations and proving theorems; it is also about developing essential cognitive skills that transcend the classroom." + `
  " The discipline of mathematics cultivates critical thinking, logical reasoning, and problem-solving abilities." + `
  " It teaches students to analyze situations systematically, identify underlying structures, and communicate ideas clearly and concisely." + `
  " These skills extend beyond mathematics, becoming valuable assets in all aspects of life, helping students navigate the complexities of adulthood"

Replace-InRange $d.Paragraphs(5).Range `
  " These particles, if they exist, could account for the observed gravitational effects of dark matter" `
  ($BRK + $newTail)

Write-Output "body-paragraph-done"

# ---------------------------------------------------------------------------
# Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------------
Replace-InRange $d.Paragraphs(7).Range `
  "In this exploration of the enigma of dark matter, the mystifying substance revealed its compelling presence through gravitational influences on visible matter" `
  "Mathematics is an elegant and intricate system of thought that captivates the human mind"

Replace-InRange $d.Paragraphs(7).Range `
  " Weighing heavily upon the universe, dark matter's existence remains concealed, its properties elusive" `
  " Its study enriches our understanding of the world and empowers us to solve problems creatively and systematically"

Replace-InRange $d.Paragraphs(7).Range `
  " The journey of uncovering its secrets continues, fueled by the ingenuity and perseverance of scientists seeking to unravel this cosmic enigma" `
  (" The aesthetic beauty of mathematics, its practical applications, and its role in developing essential thinking skills make it a subject of paramount importance." + `
   " As students embark on their mathematical journey, they discover the symphony of numbers, a harmonious composition of logic, reason, and creativity")

Write-Output "summary-paragraph-done"

# ---------------------------------------------------------------------------
# Trailing empty paragraph added at the very end of the document body.
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "done"
